# C5-PowerPoint.pptx edit
#
# 1) Table on slide 6 switches from the deck's one custom table style
#    ({43C5780E-264A-46B2-BE2C-33ACE5DBF616}) to the built-in style
#    {E400C113-A2C2-4DA8-B2BF-A90CE7AC23B6}.
#
# 2) The presentation's applied theme ("Integral", currently stored as
#    ppt/theme/theme2.xml - the theme actually wired to the slide master)
#    is switched back to the default "Office Theme" colour palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 -----------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{E400C113-A2C2-4DA8-B2BF-A90CE7AC23B6}", $false)
    }
}

# --- 2) Swap the theme colour palette back to the default Office colours
function HexToVbRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# PpColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $colorScheme.Colors($i).RGB = HexToVbRgb($officeColors[$i - 1])
}
